# Update cryptos list price (D) and 1h volume change (E) columns
# for rows 2-51 on the active worksheet, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/percentage columns keep their text representation
# (e.g. "1.00", "43.035.01") instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.035.01'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '2.301.19'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '300.60'
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").Value = '99.76'
$ws.Range("E6").Value = '  +2.51%  '

$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("D10").Value = '36.28'
$ws.Range("E10").Value = '  +7.72%  '

$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("D12").Value = '0.117'
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("E13").Value = '  +6.55%  '

$ws.Range("E14").Value = '  +2.02%  '

$ws.Range("D15").Value = '2.661.30'
$ws.Range("E15").Value = '  +0.11%  '

$ws.Range("D16").Value = '2.319.92'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("E17").Value = '  -1.35%  '

$ws.Range("D18").Value = '42.946.99'
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("E19").Value = '  +9.87%  '

$ws.Range("E20").Value = '  +0.37%  '

$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  +1.13%  '

$ws.Range("D22").Value = '67.89'
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").Value = '235.57'
$ws.Range("E23").Value = '  -0.40%  '

$ws.Range("E24").Value = '  +7.67%  '

$ws.Range("E26").Value = '  -0.86%  '

$ws.Range("D27").Value = '24.94'
$ws.Range("E27").Value = '  +1.76%  '

$ws.Range("E28").Value = '  +14.71%  '

$ws.Range("D29").Value = '34.63'
$ws.Range("E29").Value = '  +1.83%  '

$ws.Range("D30").Value = '167.63'
$ws.Range("E30").Value = '  +0.66%  '

$ws.Range("D31").Value = '9.13'
$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").Value = '5.03'
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("D34").Value = '17.61'
$ws.Range("E34").Value = '  +4.70%  '

$ws.Range("D35").Value = '4.62'
$ws.Range("E35").Value = '  -1.61%  '

$ws.Range("E36").Value = '  +1.01%  '

$ws.Range("E37").Value = '  -0.74%  '

$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("D42").Value = '2.29'
$ws.Range("E42").Value = '  -4.25%  '

$ws.Range("D43").Value = '0.0291'
$ws.Range("E43").Value = '  +3.37%  '

$ws.Range("D44").Value = '1.979.89'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("E45").Value = '  +2.95%  '

$ws.Range("E46").Value = '  +1.78%  '

$ws.Range("D47").Value = '17.48'
$ws.Range("E47").Value = '  -1.29%  '

$ws.Range("D48").Value = '55.32'
$ws.Range("E48").Value = '  +3.62%  '

$ws.Range("D49").Value = '1.55'
$ws.Range("E49").Value = '  +3.80%  '

$ws.Range("D50").Value = '2.523.89'
$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("D51").Value = '70.75'
$ws.Range("E51").Value = '  +0.93%  '
